# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to match the newly scraped snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13838
$ws1.Range("F6").Value = 497
$ws1.Range("F8").Value = 1034
$ws1.Range("F9").Value = 13903
$ws1.Range("F10").Value = 14792
$ws1.Range("F23").Value = 1148
$ws1.Range("F26").Value = 5748
$ws1.Range("F28").Value = 1058
$ws1.Range("F29").Value = 5417
$ws1.Range("F31").Value = 51
$ws1.Range("F32").Value = 273

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13838
$ws4.Range("F7").Value = 497
$ws4.Range("F9").Value = 1034
$ws4.Range("F10").Value = 13903
$ws4.Range("F11").Value = 14792
$ws4.Range("F24").Value = 1148
$ws4.Range("F27").Value = 5748
$ws4.Range("F29").Value = 1058
$ws4.Range("F30").Value = 5417
$ws4.Range("F32").Value = 51
$ws4.Range("F33").Value = 273
